$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column (D) holds values that look numeric (e.g. "578.72") but must
# stay plain text, exactly as authored in the source feed (some even contain
# two '.' separators, e.g. "69.250.34"). Force the column to Text format
# before writing so Excel doesn't silently coerce the strings to doubles.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2 - Bitcoin
$ws.Range("D2").Value = "69.250.34"
$ws.Range("E2").Value = "  -1.85%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "3.482.33"
$ws.Range("E3").Value = "  -3.77%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.01%  "

# Row 5 - BNB
$ws.Range("D5").Value = "578.72"
$ws.Range("E5").Value = "  -0.89%  "

# Row 6 - Solana
$ws.Range("D6").Value = "181.71"
$ws.Range("E6").Value = "  -4.82%  "

# Row 7 - LidoStakedEther
$ws.Range("D7").Value = "3.471.84"
$ws.Range("E7").Value = "  -4.02%  "

# Row 8 - XRP
$ws.Range("D8").Value = "0.608"
$ws.Range("E8").Value = "  -3.90%  "

# Row 9 - USDC
$ws.Range("E9").Value = "  +0.12%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  +4.81%  "

# Row 11 - Cardano
$ws.Range("D11").Value = "0.638"
$ws.Range("E11").Value = "  -3.98%  "

# Row 12 - Avalanche
$ws.Range("D12").Value = "53.49"
$ws.Range("E12").Value = "  -4.86%  "

# Row 13 - ShibaInu
$ws.Range("E13").Value = "  -3.34%  "

# Row 14 - Polkadot
$ws.Range("D14").Value = "9.35"
$ws.Range("E14").Value = "  -3.88%  "

# Row 15 - WrappedliquidstakedEther2.0
$ws.Range("D15").Value = "4.031.00"
$ws.Range("E15").Value = "  -4.01%  "

# Row 16 - Chainlink
$ws.Range("D16").Value = "19.13"
$ws.Range("E16").Value = "  -4.17%  "

# Row 17 - WrappedBTC
$ws.Range("D17").Value = "69.212.96"
$ws.Range("E17").Value = "  -1.84%  "

# Row 18 - WrappedEther
$ws.Range("D18").Value = "3.466.99"
$ws.Range("E18").Value = "  -4.09%  "

# Row 19 - Uniswap
$ws.Range("D19").Value = "12.20"
$ws.Range("E19").Value = "  -3.96%  "

# Row 20 - TRON
$ws.Range("E20").Value = "  -1.76%  "

# Row 21 - BitcoinCash
$ws.Range("D21").Value = "532.42"
$ws.Range("E21").Value = "  +7.64%  "

# Row 22 - Polygon
$ws.Range("E22").Value = "  -4.84%  "

# Row 23 - InternetComputer(DFINITY)
$ws.Range("D23").Value = "18.30"
$ws.Range("E23").Value = "  -4.77%  "

# Row 24 - PancakeSwap
$ws.Range("D24").Value = "4.47"
$ws.Range("E24").Value = "  +1.89%  "

# Row 25 - Toncoin
$ws.Range("E25").Value = "  -2.69%  "

# Row 26 - Litecoin
$ws.Range("D26").Value = "95.68"
$ws.Range("E26").Value = "  -1.88%  "

# Row 27 - ImmutableX
$ws.Range("D27").Value = "2.94"
$ws.Range("E27").Value = "  -2.05%  "

# Row 28 - RenderToken
$ws.Range("D28").Value = "10.97"
$ws.Range("E28").Value = "  -1.58%  "

# Row 29 - Filecoin
$ws.Range("D29").Value = "9.01"
$ws.Range("E29").Value = "  -4.79%  "

# Row 30 - EthereumClassic
$ws.Range("D30").Value = "31.59"
$ws.Range("E30").Value = "  -2.72%  "

# Row 31 - NEARProtocol
$ws.Range("D31").Value = "7.14"
$ws.Range("E31").Value = "  -6.07%  "

# Row 32 - Cosmos
$ws.Range("D32").Value = "12.40"
$ws.Range("E32").Value = "  +0.96%  "

# Row 33 - OKB
$ws.Range("D33").Value = "63.47"
$ws.Range("E33").Value = "  -4.13%  "

# Row 34 - Hedera
$ws.Range("E34").Value = "  -5.73%  "

# Row 35 - Bittensor
$ws.Range("D35").Value = "526.69"
$ws.Range("E35").Value = "  -9.84%  "

# Row 36 - TheGraph
$ws.Range("D36").Value = "0.404"
$ws.Range("E36").Value = "  +0.70%  "

# Row 37 - InjectiveProtocol
$ws.Range("D37").Value = "37.76"
$ws.Range("E37").Value = "  -3.60%  "

# Rows 38/39 swap - Fetch.AI and Dai exchange places, each with updated
# price and 1h volume values
$ws.Range("B38").Value = "Dai"
$ws.Range("C38").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D38").Value = "0.999"
$ws.Range("E38").Value = "  -0.26%  "

$ws.Range("B39").Value = "Fetch.AI"
$ws.Range("C39").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D39").Value = "3.04"
$ws.Range("E39").Value = "  +3.65%  "

# Row 40 - PEPE
$ws.Range("E40").Value = "  -8.20%  "

# Row 41 - Stacks
$ws.Range("E41").Value = "  -3.62%  "

# Row 42 - Maker
$ws.Range("D42").Value = "3.336.70"
$ws.Range("E42").Value = "  +2.83%  "

# Row 43 - Kaspa
$ws.Range("E43").Value = "  -2.87%  "

# Row 44 - dogwifhat
$ws.Range("D44").Value = "3.00"
$ws.Range("E44").Value = "  -7.67%  "

# Row 45 - ApeXProtocol
$ws.Range("D45").Value = "3.49"
$ws.Range("E45").Value = "  +2.28%  "

# Row 46 - ThetaToken
$ws.Range("D46").Value = "2.93"
$ws.Range("E46").Value = "  -4.69%  "

# Row 47 - VeChain
$ws.Range("D47").Value = "0.0434"
$ws.Range("E47").Value = "  -3.00%  "

# Row 48 - Stellar
$ws.Range("D48").Value = "0.133"
$ws.Range("E48").Value = "  -3.89%  "

# Row 49 - THORChain
$ws.Range("D49").Value = "8.91"
$ws.Range("E49").Value = "  -8.54%  "

# Row 50 - FirstDigitalUSD
$ws.Range("E50").Value = "  +0.06%  "

# Row 51 - Monero
$ws.Range("D51").Value = "136.10"
$ws.Range("E51").Value = "  -3.09%  "
